# Scheduled-runner update: refresh cached profit-calc outputs (cols H-N)
# on the rows whose market snapshot changed since the last run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1724.7441
$ws.Range("H18").Value = 7940.1816
$ws.Range("I18").Value = 1669.4
$ws.Range("J18").Value = 13165.833
$ws.Range("K18").Value = 1669.4
$ws.Range("L18").Value = 13165.833
$ws.Range("M18").Value = -1385.4
$ws.Range("N18").Value = -13733.833
$ws.Range("H28").Value = 771.5789
$ws.Range("I28").Value = 771.5789
$ws.Range("K28").Value = 771.5789
$ws.Range("M28").Value = -286.5789
$ws.Range("H40").Value = 2716.3333
$ws.Range("J40").Value = 2726.4
$ws.Range("L40").Value = 2726.4
$ws.Range("N40").Value = -3076.4
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H64").Value = 6787.722
$ws.Range("I64").Value = 4150.3335
$ws.Range("J64").Value = 8106.4165
$ws.Range("K64").Value = 4150.3335
$ws.Range("L64").Value = 8106.4165
$ws.Range("M64").Value = -3902.3335
$ws.Range("N64").Value = -8602.416499999999
$ws.Range("H67").Value = 6787.722
$ws.Range("I67").Value = 4150.3335
$ws.Range("J67").Value = 8106.4165
$ws.Range("K67").Value = 4150.3335
$ws.Range("L67").Value = 8106.4165
$ws.Range("M67").Value = -3292.3335
$ws.Range("N67").Value = -9822.416499999999
$ws.Range("H100").Value = 6326.125
$ws.Range("I100").Value = 2685.5715
$ws.Range("J100").Value = 9157.666999999999
$ws.Range("K100").Value = 2685.5715
$ws.Range("L100").Value = 9157.666999999999
$ws.Range("M100").Value = -2144.5715
$ws.Range("N100").Value = -10239.667
$ws.Range("H112").Value = 1622.0217
$ws.Range("J112").Value = 1740.25
$ws.Range("L112").Value = 5220.75
$ws.Range("N112").Value = -7436.75
$ws.Range("H138").Value = 3384.7693
$ws.Range("I138").Value = 2395.5293
$ws.Range("J138").Value = 4149.1816
$ws.Range("K138").Value = 7186.5879
$ws.Range("L138").Value = 12447.5448
$ws.Range("M138").Value = -2046.5879
$ws.Range("N138").Value = -22727.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5749324.5
$ws.Range("I32").Value = 5884546
$ws.Range("J32").Value = 2402
$ws.Range("K32").Value = 5884546
$ws.Range("L32").Value = 2402
$ws.Range("M32").Value = -5884259
$ws.Range("N32").Value = -2976
$ws.Range("H122").Value = 4931.8335
$ws.Range("I122").Value = 4928.7856
$ws.Range("K122").Value = 14786.3568
$ws.Range("M122").Value = -12336.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1958.6666
$ws.Range("I20").Value = 1884.9
$ws.Range("J20").Value = 2050.875
$ws.Range("K20").Value = 1884.9
$ws.Range("L20").Value = 2050.875
$ws.Range("M20").Value = -1637.9
$ws.Range("N20").Value = -2544.875
$ws.Range("H105").Value = 1576.9259
$ws.Range("I105").Value = 1611.5416
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1611.5416
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 135.4584
$ws.Range("N105").Value = -4794
$ws.Range("H107").Value = 6129.5
$ws.Range("J107").Value = 6499.4
$ws.Range("L107").Value = 6499.4
$ws.Range("N107").Value = -10339.4
$ws.Range("H134").Value = 423522.06
$ws.Range("I134").Value = 518261.97
$ws.Range("K134").Value = 1554785.91
$ws.Range("M134").Value = -1552250.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 72353.06
$ws.Range("J68").Value = 72353.06
$ws.Range("L68").Value = 72353.06
$ws.Range("N68").Value = -73851.06
$ws.Range("H71").Value = 72353.06
$ws.Range("J71").Value = 72353.06
$ws.Range("L71").Value = 217059.18
$ws.Range("N71").Value = -224547.18
$ws.Range("H107").Value = 1052.625
$ws.Range("I107").Value = 1052.625
$ws.Range("K107").Value = 1052.625
$ws.Range("M107").Value = 867.375
$ws.Range("H119").Value = 98340.664
$ws.Range("J119").Value = 98340.664
$ws.Range("L119").Value = 98340.664
$ws.Range("N119").Value = -108016.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5297
$ws.Range("I80").Value = 4995
$ws.Range("J80").Value = 5448
$ws.Range("K80").Value = 14985
$ws.Range("L80").Value = 16344
$ws.Range("M80").Value = -14049
$ws.Range("N80").Value = -18216
$ws.Range("H83").Value = 5297
$ws.Range("I83").Value = 4995
$ws.Range("J83").Value = 5448
$ws.Range("K83").Value = 44955
$ws.Range("L83").Value = 49032
$ws.Range("M83").Value = -40275
$ws.Range("N83").Value = -58392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1512.75
$ws.Range("I102").Value = 1086.6111
$ws.Range("K102").Value = 1086.6111
$ws.Range("M102").Value = 535.3888999999999
$ws.Range("H132").Value = 755062.6
$ws.Range("J132").Value = 3657
$ws.Range("L132").Value = 10971
$ws.Range("N132").Value = -16031

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2740.25
$ws.Range("I7").Value = 2703.1428
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 2703.1428
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -2591.1428
$ws.Range("N7").Value = -3224
$ws.Range("H55").Value = 1198.1
$ws.Range("I55").Value = 176.42857
$ws.Range("K55").Value = 176.42857
$ws.Range("M55").Value = -3.428570000000008
$ws.Range("H68").Value = 2889.05
$ws.Range("I68").Value = 2652.6924
$ws.Range("J68").Value = 3328
$ws.Range("K68").Value = 2652.6924
$ws.Range("L68").Value = 3328
$ws.Range("M68").Value = -1903.6924
$ws.Range("N68").Value = -4826
$ws.Range("H71").Value = 2889.05
$ws.Range("I71").Value = 2652.6924
$ws.Range("J71").Value = 3328
$ws.Range("K71").Value = 13263.462
$ws.Range("L71").Value = 16640
$ws.Range("M71").Value = -9519.462
$ws.Range("N71").Value = -24128
$ws.Range("H122").Value = 3570.8408
$ws.Range("I122").Value = 3345.4285
$ws.Range("K122").Value = 10036.2855
$ws.Range("M122").Value = -7586.2855
$ws.Range("H126").Value = 2740.25
$ws.Range("I126").Value = 2703.1428
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8109.428400000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5639.428400000001
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1734.0857
$ws.Range("I122").Value = 1291.3846
$ws.Range("J122").Value = 3013
$ws.Range("K122").Value = 3874.1538
$ws.Range("L122").Value = 9039
$ws.Range("M122").Value = -1424.1538
$ws.Range("N122").Value = -13939
$ws.Range("H126").Value = 4673.625
$ws.Range("I126").Value = 3912.8572
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 11738.5716
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -9268.571599999999
$ws.Range("N126").Value = -34937
$ws.Range("H136").Value = 10561754
$ws.Range("I136").Value = 12672939
$ws.Range("K136").Value = 38018817
$ws.Range("M136").Value = -38016267
